$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy style from an existing header cell (e.g. AB1) to the new header cells
$ws.Range("AB1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data rows 2 to 41: Wins = 93, Losses = 69, Ties = 0
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 29).Value = 93   # AC
    $ws.Cells.Item($r, 30).Value = 69   # AD
    $ws.Cells.Item($r, 31).Value = 0    # AE
}
